$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.630.02"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = "'2.363.08"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.52%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').Value = "'0.662"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.37%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').Value = "'234.23"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.65%  '
$ws.Range('D7').Value = "'73.41"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +13.85%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = "'0.518"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +18.57%  '
$ws.Range('D10').Value = "'0.0979"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.28%  '
$ws.Range('D11').Value = "'27.53"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.80%  '
$ws.Range('D12').Value = "'2.723.85"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.04%  '
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('D14').Value = "'16.42"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.29%  '
$ws.Range('D15').Value = "'6.46"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.69%  '
$ws.Range('D16').Value = "'0.879"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.34%  '
$ws.Range('D17').Value = "'2.367.14"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.68%  '
$ws.Range('D18').Value = "'43.566.07"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('D19').Value = "'0.0000100"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.22%  '
$ws.Range('D20').Value = "'75.60"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.87%  '
$ws.Range('D21').Value = "'6.39"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.97%  '
$ws.Range('D22').Value = "'251.64"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.53%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = "'1.00"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('B24').Value = 'WEMIXToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D24').Value = "'3.76"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('E25').Value = '  +2.98%  '
$ws.Range('D26').Value = "'10.18"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.60%  '
$ws.Range('D27').Value = "'2.24"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.40%  '
$ws.Range('D28').Value = "'22.54"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.51%  '
$ws.Range('D29').Value = "'172.02"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('D30').Value = "'1.53"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +9.01%  '
$ws.Range('D31').Value = "'0.134"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.58%  '
$ws.Range('D32').Value = "'0.130"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.51%  '
$ws.Range('D33').Value = "'5.11"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.32%  '
$ws.Range('D34').Value = "'0.0697"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.20%  '
$ws.Range('D35').Value = "'5.09"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.35%  '
$ws.Range('D36').Value = "'3.78"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.42%  '
$ws.Range('D37').Value = "'6.66"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.15%  '
$ws.Range('D38').Value = "'2.44"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.47%  '
$ws.Range('E39').Value = '  +4.90%  '
$ws.Range('D40').Value = "'19.63"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +15.94%  '
$ws.Range('B41').Value = 'BinanceUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D41').Value = "'1.00"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = "'8.86"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.67%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = "'1.17"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +11.00%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = "'100.06"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.77%  '
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('E46').Value = '  +3.45%  '
$ws.Range('E47').Value = '  +2.77%  '
$ws.Range('D48').Value = "'1.445.49"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('E49').Value = '  +9.98%  '
$ws.Range('D50').Value = "'2.593.38"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.06%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').Value = "'2.77"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.07%  '
